$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 used to hold a stray pandas multi-index artifact ("unnamed: 1_level_1");
# correct it to "total" to match column B1 above it.
$ws.Range("B2").Value = "total"

# Delete the "situação do domicílio" label-only row (original row 5)
$ws.Rows.Item(5).Delete()

# After the above deletion, "grandes regiões e unidades da federação" shifted
# up from row 8 to row 7. Delete it too.
$ws.Rows.Item(7).Delete()
